# Fix: Springer API author strings for E2:E14 were regenerated with corrected
# whitespace padding between author entries (bug in abstract-fill logic fixed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Rong%Chen%NULL%2,          Yuan%Zhang%NULL%1,          Lei%Huang%NULL%0,          Bi-heng%Cheng%NULL%1,          Zhong-yuan%Xia%NULL%1,          Qing-tao%Meng%mengqingtao2018@126.com%1]'
$ws.Range("E3").Value = '[Lingkong%Zeng%xref no email%2,    Shiwen%Xia%xref no email%3,    Wenhao%Yuan%xref no email%2,    Kai%Yan%xref no email%2,    Feifan%Xiao%xref no email%2,    Jianbo%Shao%xref no email%2,    Wenhao%Zhou%xref no email%3]'
$ws.Range("E4").Value = '[S.%Khan%Suliman.khan18@mails.ucas.ac.cn%1,          L.%Jun%NULL%1,          NULL%Nawsherwan%NULL%1,          R.%Siddique%NULL%1,          Y.%Li%NULL%1,          G.%Han%hg7913@hotmail.com%1,          M.%Xue%xuemengzhou@zzu.edu.cn%1,          G.%Nabi%ghulamnabiqau@gmail.com%1,          J.%Liu%jbliuzz@163.com%1]'
$ws.Range("E5").Value = '[Xiaoqing%Wu%NULL%1,          Ruihong%Sun%NULL%1,          Jianpu%Chen%NULL%1,          Yuanliang%Xie%NULL%1,          Shutong%Zhang%NULL%1,          Xiang%Wang%wangxiang_897@126.com%1]'
$ws.Range("E6").Value = '[Lian%Chen%NULL%1,          Qin%Li%NULL%1,          Danni%Zheng%NULL%1,          Hai%Jiang%NULL%1,          Yuan%Wei%NULL%0,          Li%Zou%NULL%2,          Li%Zou%NULL%0,          Ling%Feng%NULL%3,          Ling%Feng%NULL%0,          Guoping%Xiong%NULL%3,          Guoping%Xiong%NULL%0,          Guoqiang%Sun%NULL%4,          Guoqiang%Sun%NULL%0,          Haibo%Wang%NULL%2,          Haibo%Wang%NULL%0,          Yangyu%Zhao%NULL%2,          Yangyu%Zhao%NULL%0,          Jie%Qiao%NULL%1]'
$ws.Range("E7").Value = '[Jie%Yan%NULL%2,          Juanjuan%Guo%NULL%0,          Cuifang%Fan%NULL%1,          Juan%Juan%NULL%1,          Xuechen%Yu%NULL%0,          Jiafu%Li%NULL%0,          Ling%Feng%NULL%0,          Chunyan%Li%NULL%1,          Huijun%Chen%NULL%1,          Yuan%Qiao%NULL%1,          Di%Lei%NULL%1,          Chen%Wang%NULL%0,          Guoping%Xiong%NULL%0,          Fengyi%Xiao%NULL%1,          Wencong%He%NULL%1,          Qiumei%Pang%NULL%1,          Xiaoling%Hu%NULL%1,          Suqing%Wang%NULL%1,          Dunjin%Chen%NULL%1,          Yuanzhen%Zhang%NULL%0,          Liona C.%Poon%NULL%1,          Huixia%Yang%NULL%0]'
$ws.Range("E8").Value = '[Dehan%Liu%xref no email%1,    Lin%Li%xref no email%1,    Xin%Wu%xref no email%1,    Dandan%Zheng%xref no email%1,    Jiazheng%Wang%xref no email%1,    Lian%Yang%xref no email%1,    Chuansheng%Zheng%xref no email%1]'
$ws.Range("E9").Value = '[Na%Li%NULL%4,          Lefei%Han%NULL%1,          Min%Peng%2658706528@qq.com%2,          Yuxia%Lv%NULL%1,          Yin%Ouyang%NULL%1,          Kui%Liu%NULL%0,          Linli%Yue%NULL%1,          Qiannan%Li%NULL%1,          Guoqiang%Sun%NULL%0,          Lin%Chen%NULL%1,          Lin%Yang%l.yang@polyu.edu.hk%2]'
$ws.Range("E10").Value = '[ Y.%Zhang%null%1,           R.% Chen%null%1,           J.% Wang%null%1,           Y.% Gong%null%1,           Q.% Zhou%null%1,           H.-h.% Cheng%null%1,           Z.-y.% Xia%null%1,           X.% Chen%null%1,           Q.-t.% Meng%null%1,           D. % Ma%null%1]'
$ws.Range("E11").Value = '[Enrico M.%Ferrazzi%enrico.ferrazzi@unimi.it%2,          Luigi%Frigerio%NULL%1,          Irene%Cetin%NULL%1,          Patrizia%Vergani%NULL%1,          Arsenio%Spinillo%NULL%1,          Federico%Prefumo%NULL%1,          Edda%Pellegrini%NULL%1,          Gianluigi%Gargantini%NULL%1]'
$ws.Range("E12").Value = '[Noelle%Breslin%nb2565@cumc.columbia.edu%3,          Caitlin%Baptiste%NULL%1,          Cynthia%Gyamfi-Bannerman%NULL%1,          Russell%Miller%NULL%1,          Rebecca%Martinez%NULL%1,          Kyra%Bernstein%NULL%1,          Laurence%Ring%NULL%1,          Ruth%Landau%NULL%1,          Stephanie%Purisch%NULL%1,          Alexander M.%Friedman%NULL%1,          Karin%Fuchs%NULL%2,          Desmond%Sutton%NULL%1,          Maria%Andrikopoulou%NULL%1,          Devon%Rupley%NULL%1,          Jean-Ju%Sheen%NULL%1,          Janice%Aubey%NULL%1,          Noelia%Zork%NULL%1,          Leslie%Moroz%NULL%1,          Mirella%Mourad%NULL%1,          Ronald%Wapner%NULL%1,          Lynn L.%Simpson%NULL%1,          Mary E.%D’Alton%NULL%1,          Dena%Goffman%NULL%2]'
$ws.Range("E13").Value = '[Desmond%Sutton%NULL%1,          Karin%Fuchs%NULL%0,          Mary%D’Alton%NULL%1,          Dena%Goffman%NULL%0]'
$ws.Range("E14").Value = '[Wei%Liu%NULL%0,          Jing%Wang%NULL%1,          Wenbin%Li%NULL%1,          Zhaoxian%Zhou%NULL%1,          Siying%Liu%NULL%1,          Zhihui%Rong%rongzhihui53@163.com%1]'
